# Rename sector labels in column D ("sector") to their updated names.
#   Communication -> Communication Services  (rows 2-21)
#   Health Care   -> Healthcare              (rows 22-37, 142-150)
#   Materials     -> Basic Materials         (rows 38-106)
# Rows 107-141 (Industrials) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = "Communication Services"
}

for ($r = 22; $r -le 37; $r++) {
    $ws.Cells.Item($r, 4).Value = "Healthcare"
}

for ($r = 38; $r -le 106; $r++) {
    $ws.Cells.Item($r, 4).Value = "Basic Materials"
}

for ($r = 142; $r -le 150; $r++) {
    $ws.Cells.Item($r, 4).Value = "Healthcare"
}
